$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Case_3_175 (380 kV case): refreshed bus voltage magnitudes (res_bus/vm_pu)
# Slack bus setpoint (column B) lowered from 1.05 to 1.02 p.u.; all other bus
# voltages recomputed accordingly for rows 2-25 (bus indices 0-23).

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028410913570541
$ws.Range("D2").Value = 1.031890569828557
$ws.Range("E2").Value = 1.037913677981729
$ws.Range("F2").Value = 1.047886746308251
$ws.Range("I2").Value = 1.033590634667089
$ws.Range("J2").Value = 1.033563382985585
$ws.Range("K2").Value = 1.034697380732762
$ws.Range("L2").Value = 1.040703205692687
$ws.Range("M2").Value = 1.050648130421237
$ws.Range("N2").Value = 1.015205161976223

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029182016598456
$ws.Range("D3").Value = 1.032440308949029
$ws.Range("E3").Value = 1.038673012206232
$ws.Range("F3").Value = 1.048895550491063
$ws.Range("I3").Value = 1.033726371819864
$ws.Range("J3").Value = 1.03397599804757
$ws.Range("K3").Value = 1.035056471724849
$ws.Range("L3").Value = 1.041272579211842
$ws.Range("M3").Value = 1.051468359076404
$ws.Range("N3").Value = 1.015342669699636

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029681573192611
$ws.Range("D4").Value = 1.032796485322546
$ws.Range("E4").Value = 1.039165345577373
$ws.Range("F4").Value = 1.049549857184679
$ws.Range("I4").Value = 1.033813274067905
$ws.Range("J4").Value = 1.034242907764914
$ws.Range("K4").Value = 1.035288578467382
$ws.Range("L4").Value = 1.041641331708275
$ws.Range("M4").Value = 1.052000031805813
$ws.Range("N4").Value = 1.015431593668862

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029891728996561
$ws.Range("D5").Value = 1.032946329838063
$ws.Range("E5").Value = 1.039372558535341
$ws.Range("F5").Value = 1.049825295299053
$ws.Range("I5").Value = 1.033849584756563
$ws.Range("J5").Value = 1.034355096292545
$ws.Range("K5").Value = 1.03538609529812
$ws.Range("L5").Value = 1.041796432791197
$ws.Range("M5").Value = 1.052223768014074
$ws.Range("N5").Value = 1.015468964246566

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02992702336907
$ws.Range("D6").Value = 1.03297149567121
$ws.Range("E6").Value = 1.039407364290051
$ws.Range("F6").Value = 1.049871564073288
$ws.Range("I6").Value = 1.033855668383786
$ws.Range("J6").Value = 1.034373932013416
$ws.Range("K6").Value = 1.035402465200318
$ws.Range("L6").Value = 1.041822479451094
$ws.Range("M6").Value = 1.052261347219525
$ws.Range("N6").Value = 1.015475238152505

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029684380748154
$ws.Range("D7").Value = 1.032798487130291
$ws.Range("E7").Value = 1.039168113441967
$ws.Range("F7").Value = 1.049553536160154
$ws.Range("I7").Value = 1.033813760130238
$ws.Range("J7").Value = 1.034244406914861
$ws.Range("K7").Value = 1.035289881732626
$ws.Range("L7").Value = 1.041643403873651
$ws.Range("M7").Value = 1.052003020513675
$ws.Range("N7").Value = 1.015432093068212

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028671385873204
$ws.Range("D8").Value = 1.032076260756982
$ws.Range("E8").Value = 1.038170091774345
$ws.Range("F8").Value = 1.048227355847264
$ws.Range("I8").Value = 1.033636699578599
$ws.Range("J8").Value = 1.033702843977756
$ws.Range("K8").Value = 1.034818788080378
$ws.Range("L8").Value = 1.040895559012319
$ws.Range("M8").Value = 1.050925136999394
$ws.Range("N8").Value = 1.01525164400906

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026891043116869
$ws.Range("D9").Value = 1.030807194825246
$ws.Range("E9").Value = 1.036419137713559
$ws.Range("F9").Value = 1.045902353495902
$ws.Range("I9").Value = 1.033317612621942
$ws.Range("J9").Value = 1.032747983052324
$ws.Range("K9").Value = 1.033986805775386
$ws.Range("L9").Value = 1.039580348816103
$ws.Range("M9").Value = 1.049032963323709
$ws.Range("N9").Value = 1.014933285240007

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025707401980959
$ws.Range("D10").Value = 1.029963670059794
$ws.Range("E10").Value = 1.035257111423203
$ws.Range("F10").Value = 1.044360458209668
$ws.Range("I10").Value = 1.033100165342757
$ws.Range("J10").Value = 1.032111106256426
$ws.Range("K10").Value = 1.033430975278765
$ws.Range("L10").Value = 1.038705366352384
$ws.Range("M10").Value = 1.047776449404648
$ws.Range("N10").Value = 1.014720813654076

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025195666286411
$ws.Range("D11").Value = 1.029599034507575
$ws.Range("E11").Value = 1.034755215485144
$ws.Range("F11").Value = 1.043694745900673
$ws.Range("I11").Value = 1.033004896573398
$ws.Range("J11").Value = 1.03183527500429
$ws.Range("K11").Value = 1.033190030847708
$ws.Range("L11").Value = 1.038326939806344
$ws.Range("M11").Value = 1.047233556204585
$ws.Range("N11").Value = 1.014628761502957

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02500570515121
$ws.Range("D12").Value = 1.029463687039066
$ws.Range("E12").Value = 1.034568981506605
$ws.Range("F12").Value = 1.043447763782221
$ws.Range("I12").Value = 1.032969343015447
$ws.Range("J12").Value = 1.031732811302365
$ws.Range("K12").Value = 1.033100494592498
$ws.Range("L12").Value = 1.038186443877881
$ws.Range("M12").Value = 1.04703208128909
$ws.Range("N12").Value = 1.014594562085234

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025046446987568
$ws.Range("D13").Value = 1.0294927152289
$ws.Range("E13").Value = 1.034608920585563
$ws.Range("F13").Value = 1.043500728993345
$ws.Range("I13").Value = 1.032976976897805
$ws.Range("J13").Value = 1.031754790442645
$ws.Range("K13").Value = 1.033119702157788
$ws.Range("L13").Value = 1.038216577606769
$ws.Range("M13").Value = 1.047075290179672
$ws.Range("N13").Value = 1.014601898292852

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025179961575546
$ws.Range("D14").Value = 1.029587844707092
$ws.Range("E14").Value = 1.034739817388339
$ws.Range("F14").Value = 1.04367432429248
$ws.Range("I14").Value = 1.033001961097778
$ws.Range("J14").Value = 1.031826805472555
$ws.Range("K14").Value = 1.033182630536662
$ws.Range("L14").Value = 1.038315324962549
$ws.Range("M14").Value = 1.047216898547886
$ws.Range("N14").Value = 1.01462593471171

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025262240312463
$ws.Range("D15").Value = 1.029646469683263
$ws.Range("E15").Value = 1.034820492794737
$ws.Range("F15").Value = 1.043781320975683
$ws.Range("I15").Value = 1.033017332645801
$ws.Range("J15").Value = 1.031871175326614
$ws.Range("K15").Value = 1.033221397687518
$ws.Range("L15").Value = 1.038376175592798
$ws.Range("M15").Value = 1.047304172011184
$ws.Range("N15").Value = 1.0146407434085

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025741380963925
$ws.Range("D16").Value = 1.029987882860063
$ws.Range("E16").Value = 1.035290447500322
$ws.Range("F16").Value = 1.04440468034195
$ws.Range("I16").Value = 1.033106464631837
$ws.Range("J16").Value = 1.032129411121788
$ws.Range("K16").Value = 1.033446960463784
$ws.Range("L16").Value = 1.038730490808761
$ws.Range("M16").Value = 1.047812504560069
$ws.Range("N16").Value = 1.014726921825072

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026042145663222
$ws.Range("D17").Value = 1.030202208602003
$ws.Range("E17").Value = 1.0355855787547
$ws.Range("F17").Value = 1.044796217493789
$ws.Range("I17").Value = 1.033162077352829
$ws.Range("J17").Value = 1.032291380579828
$ws.Range("K17").Value = 1.033588379638667
$ws.Range("L17").Value = 1.03895286397472
$ws.Range("M17").Value = 1.048131686736438
$ws.Range("N17").Value = 1.014780966048763

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026217652604506
$ws.Range("D18").Value = 1.03032728051677
$ws.Range("E18").Value = 1.035757846149182
$ws.Range("F18").Value = 1.045024781406561
$ws.Range("I18").Value = 1.033194407911288
$ws.Range("J18").Value = 1.032385848810967
$ws.Range("K18").Value = 1.033670841288454
$ws.Range("L18").Value = 1.039082613477426
$ws.Range("M18").Value = 1.048317974582888
$ws.Range("N18").Value = 1.014812484222715

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026277508769221
$ws.Range("D19").Value = 1.030369936816227
$ws.Range("E19").Value = 1.035816605571425
$ws.Range("F19").Value = 1.045102747442779
$ws.Range("I19").Value = 1.033205413558831
$ws.Range("J19").Value = 1.032418059013023
$ws.Range("K19").Value = 1.033698954173119
$ws.Range("L19").Value = 1.039126861954159
$ws.Range("M19").Value = 1.0483815132582
$ws.Range("N19").Value = 1.014823230259619

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02600986858545
$ws.Range("D20").Value = 1.030179207327704
$ws.Range("E20").Value = 1.035553901298
$ws.Range("F20").Value = 1.04475418991601
$ws.Range("I20").Value = 1.033156121738322
$ws.Range("J20").Value = 1.032274003382281
$ws.Range("K20").Value = 1.033573209355483
$ws.Range("L20").Value = 1.038929000980749
$ws.Range("M20").Value = 1.048097429676772
$ws.Range("N20").Value = 1.014775168119415

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025140641522267
$ws.Range("D21").Value = 1.029559828825937
$ws.Range("E21").Value = 1.034701266213095
$ws.Range("F21").Value = 1.043623196698619
$ws.Range("I21").Value = 1.032994608464296
$ws.Range("J21").Value = 1.031805599043673
$ws.Range("K21").Value = 1.03316410075784
$ws.Range("L21").Value = 1.038286244423885
$ws.Range("M21").Value = 1.047175193436429
$ws.Range("N21").Value = 1.014618856779258

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02459482079765
$ws.Range("D22").Value = 1.029170948113745
$ws.Range("E22").Value = 1.034166295469926
$ws.Range("F22").Value = 1.042913792919853
$ws.Range("I22").Value = 1.032892095701424
$ws.Range("J22").Value = 1.031511050741843
$ws.Range("K22").Value = 1.032906654264285
$ws.Range("L22").Value = 1.037882515130696
$ws.Range("M22").Value = 1.046596387952868
$ws.Range("N22").Value = 1.01452053648689

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024884104041606
$ws.Range("D23").Value = 1.029377048665636
$ws.Range("E23").Value = 1.034449787297332
$ws.Range("F23").Value = 1.043289699859253
$ws.Range("I23").Value = 1.032946530709773
$ws.Range("J23").Value = 1.031667200120416
$ws.Range("K23").Value = 1.033043152316129
$ws.Range("L23").Value = 1.038096501531139
$ws.Range("M23").Value = 1.046903124455032
$ws.Range("N23").Value = 1.014572661686183

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026024452973005
$ws.Range("D24").Value = 1.030189600430389
$ws.Range("E24").Value = 1.035568214599413
$ws.Range("F24").Value = 1.044773179792034
$ws.Range("I24").Value = 1.033158813156098
$ws.Range("J24").Value = 1.03228185540798
$ws.Range("K24").Value = 1.033580064233216
$ws.Range("L24").Value = 1.038939783507882
$ws.Range("M24").Value = 1.048112908614545
$ws.Range("N24").Value = 1.014777787969338

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027350738545401
$ws.Range("D25").Value = 1.031134842511985
$ws.Range("E25").Value = 1.036870879189538
$ws.Range("F25").Value = 1.046502001437841
$ws.Range("I25").Value = 1.033400939581163
$ws.Range("J25").Value = 1.032994896155773
$ws.Range("K25").Value = 1.034202105599257
$ws.Range("L25").Value = 1.039920046803025
$ws.Range("M25").Value = 1.049521272887373
$ws.Range("N25").Value = 1.015015631485208
